$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 (B10/C10): "Estudo de Casos..." -> "5840521 - Rosa Ana Conte"
$ws.Range("B10").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C10").Value = "5840521 - Rosa Ana Conte"

# Row 13 (B13/C13): "60 h" -> "Semestral"
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 (B15/C15): "Semestral" -> "01/01/2023"
# (Assigning the literal text "01/01/2023" via .Value would make Excel coerce
# it into a real date serial number and change the cell style. Instead copy
# the existing text value "01/01/2023" already stored in B8/C8 and paste only
# the value, which keeps it as plain text and leaves formatting untouched.)
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 (B18/C18): "01/01/2023" -> "5840521 - Rosa Ana Conte"
$ws.Range("B18").Value = "5840521 - Rosa Ana Conte"
$ws.Range("C18").Value = "5840521 - Rosa Ana Conte"
